# Vietnamese translation pass for the Partner "reminder to RSVP" email template.
$d = $word.ActiveDocument

# Replace the first occurrence of $old inside $range with $new. We locate
# the match with Find (no built-in replace) and then assign the new text
# straight onto the matched sub-range in a single step. A single Range.Text
# assignment is what keeps structural markers (<w:commentRangeStart/>,
# <w:commentRangeEnd/>, <w:commentReference/>, ...) that sit exactly at a
# match boundary anchored in their original place, regardless of whether
# the boundary is at the start or the end of the match.
function Replace-InRange($range, $old, $new) {
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $range.Text = $new
}

function Replace-InParagraph($index, $old, $new) {
    $p = $d.Paragraphs.Item($index)
    Replace-InRange $p.Range $old $new
}

# Para 12: "Subject line: Reminder: RSVP for [EVENT NAME]  "
Replace-InParagraph 12 ": Reminder: RSVP for " ": Nhắc nhở: Xác nhận tham dự sự kiện "

# Para 15: "Don't delay! Book your spot today!"
Replace-InParagraph 15 "Don’t delay! Book your spot today!" "Đừng chậm trễ! Hãy đặt chỗ của bạn ngay hôm nay!"

# Para 17: "Hi [PARTNER NAME], "
Replace-InParagraph 17 "Hi " "Xin chào "

# Para 19: "We hope you're as excited as us for the [EVENT NAME], happening on [DD Mmm YYYY]!"
Replace-InParagraph 19 "We hope you’re as excited as us for the " "Chúng tôi hy vọng bạn mong chờ sự kiện "
Replace-InParagraph 19 ", happening on " " diễn ra vào ngày "
Replace-InParagraph 19 "!" " như chúng tôi!"

# Para 20: "We hope you're as excited as us for the [EVENT NAME], happening from [DD Mmm YYYY] to [DD Mmm YYYY]!"
Replace-InParagraph 20 "We hope you’re as excited as us for the " "Chúng tôi hy vọng bạn mong chờ sự kiện "
Replace-InParagraph 20 ", happening from " " diễn ra từ ngày "
Replace-InParagraph 20 " to " " đến ngày "
Replace-InParagraph 20 "!" " như chúng tôi!"

# Para 21: "Confirm your attendance for this highly-anticipated event by [DD Mmm YYYY] as spots are limited and on a first-come, first-served basis."
Replace-InParagraph 21 "Confirm your attendance for this highly-anticipated event by [" "Hãy xác nhận tham dự sự kiện rất được mong đợi này trước ngày ["
Replace-InParagraph 21 "] as spots are limited and on a first-come, first-served basis." "] vì số lượng tham dự là có hạn và chúng tôi sẽ chọn những đơn đăng ký sớm nhất."

# Para 22: "RVSP now"
Replace-InParagraph 22 "RVSP now" "Xác nhận tham dự ngay bây giờ"

# Para 25: "If you have any questions, please contact us via live chat or WhatsApp. "
Replace-InParagraph 25 "If you have any questions, please contact us via " "Nếu bạn cần hỗ trợ, hãy liên hệ với chúng tôi qua "
Replace-InParagraph 25 " or " " hoặc "

# Para 26: "If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). "
Replace-InParagraph 26 "If you have any questions, please contact your country manager, " "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn "
Replace-InParagraph 26 ", at " ", qua email "
Replace-InParagraph 26 " or " " hoặc số "

# Para 27: "We look forward to seeing you at [EVENT NAME]! "
Replace-InParagraph 27 "We look forward to seeing you at " "Chúng tôi rất mong được gặp bạn tại sự kiện "

# Para 38: "We look forward to seeing you at [EVENT NAME]! "
Replace-InParagraph 38 "We look forward to seeing you at [EVENT NAME]! " "Chúng tôi rất mong được gặp bạn tại sự kiện [EVENT NAME]! "

# Para 39: "If you have any questions, please contact your country manager:"
Replace-InParagraph 39 "If you have any questions, please contact your country manager:" "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn:"

# Para 41: "If you have any questions, please contact us via:"
Replace-InParagraph 41 "If you have any questions, please contact us via:" "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua:"

# Para 42: "live chat | WhatsApp." -- capitalize only the "live chat" field-result text
Replace-InParagraph 42 "live chat" "Live chat"
